$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (TP2): add the comment about the redone LCOM1 presentation,
#     give it a recalculated average formula, and flag it "*Reavaliado" ---
# (values are written in shared-string order: "*Reavaliado" first, then the
# long comment, so the two new strings land at the same table slots the
# original author's save produced)
$ws.Range("D6").Value = "*Reavaliado"

$lcom1Text = @"
Implementação do LCOM1 (os dois apresentaram)
- Explicação do LCOM1 podia ser mais rápida, já que foi a única métrica que expliquei em sala de aula.
- Zona de conforto, fez a única métrica que expliquei.
- Bom exemplo. Mostrou que entendeu.
- Usou listener
"@
$ws.Range("B6").Value = $lcom1Text

# Formatting is applied in cellXf order: the wrapped comment style first,
# then the colored flag style.
$ws.Range("B6").WrapText = $true
# Highlight the "*Reavaliado" flag in the orange "Accent 2, Lighter 40%"
# font color (theme accent2 + ~40% tint) used for re-evaluation remarks.
$ws.Range("D6").Font.Color = 0x83B1F4

$ws.Range("C6").Formula = "=(0+9)/2"

# --- Row 7 (TP3): fill in the Move Method comments and the grade ---
$moveMethodText = @"
Move Method para feature envy
- Zona de conforto, foco no que eu falei em sala de aula.
- Explicação "bacana" de move method.
- Como sugeriu o move method ficou meio fraca.
- Usou listener
- ferramenta fraca, podia ser bem melhor e o exemplo podia ser maior (em termos de código).
- Legal apontar as dificuldades.
"@
$ws.Range("B7").Value = $moveMethodText
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value = 8

# Row 7 grew taller to fit the new comment text.
$ws.Rows.Item(7).RowHeight = 144

# Leave the selection where the author left it.
$ws.Range("C7").Select()
